# Updated cryptos list on Tue May 28 08:40:57 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every touched cell to remain text, matching the source feeds string cells
# (avoids Excel auto-converting numeric-looking strings like "0.998" or "598.43" to numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.954.29'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.863.93'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.82%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.43'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.35'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.861.32'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.80%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.39'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000247'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.89'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.508.82'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.863.33'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.69%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.989.34'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.13'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +5.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.37'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.95%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.90'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '465.06'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.727'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.54%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.20'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.25'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.09'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.60%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.80%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.009.21'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.94%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.31'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.13'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.35'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.79%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.835.23'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.32%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.33%  '

# Rows 38 and 39 swap coin identity (dogwifhat moves above Mantle)
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.38'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.15%  '

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.02'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.140'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.23%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '428.00'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.35%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.66%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '47.20'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.67%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000275'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '40.52'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.82'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.46%  '
